# Scheduled-runner data refresh for Sheets/Ultros_Profits.xlsx
# Updates per-item market-board price/profit columns (H:N) on each of the
# 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 358
$ws.Range("I9").Value = 258.9
$ws.Range("J9").Value = 499.57144
$ws.Range("K9").Value = 258.9
$ws.Range("L9").Value = 499.57144
$ws.Range("M9").Value = -89.89999999999998
$ws.Range("N9").Value = -837.5714399999999

$ws.Range("H28").Value = 779.63635
$ws.Range("I28").Value = 630
$ws.Range("K28").Value = 630
$ws.Range("M28").Value = -145

$ws.Range("H107").Value = 1421.9474
$ws.Range("J107").Value = 559.8333
$ws.Range("L107").Value = 559.8333
$ws.Range("N107").Value = -4399.8333

$ws.Range("H132").Value = 14141.512
$ws.Range("I132").Value = 1772.2188
$ws.Range("K132").Value = 5316.6564
$ws.Range("M132").Value = -2786.6564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4974.392
$ws.Range("I32").Value = 5086.69
$ws.Range("K32").Value = 5086.69
$ws.Range("M32").Value = -4799.69

$ws.Range("H61").Value = 1407
$ws.Range("I61").Value = 1508.1111
$ws.Range("J61").Value = 952
$ws.Range("K61").Value = 1508.1111
$ws.Range("L61").Value = 952
$ws.Range("M61").Value = -1296.1111
$ws.Range("N61").Value = -1376

$ws.Range("H74").Value = 2241.6316
$ws.Range("I74").Value = 2311.9375
$ws.Range("K74").Value = 2311.9375
$ws.Range("M74").Value = -1437.9375

$ws.Range("H77").Value = 2241.6316
$ws.Range("I77").Value = 2311.9375
$ws.Range("K77").Value = 11559.6875
$ws.Range("M77").Value = -7191.6875

$ws.Range("H136").Value = 1407
$ws.Range("I136").Value = 1508.1111
$ws.Range("J136").Value = 952
$ws.Range("K136").Value = 4524.3333
$ws.Range("L136").Value = 2856
$ws.Range("M136").Value = -1974.3333
$ws.Range("N136").Value = -7956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 82986.664
$ws.Range("J81").Value = 82986.664
$ws.Range("L81").Value = 82986.664
$ws.Range("N81").Value = -85108.664

$ws.Range("H84").Value = 82986.664
$ws.Range("J84").Value = 82986.664
$ws.Range("L84").Value = 248959.992
$ws.Range("N84").Value = -259567.992

$ws.Range("H99").Value = 102977.5
$ws.Range("I99").Value = 102977.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 102977.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -101479.5
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 2452.5
$ws.Range("I105").Value = 2270
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2270
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -523
$ws.Range("N105").Value = -6494

$ws.Range("H107").Value = 4624.185
$ws.Range("I107").Value = 4128.7827
$ws.Range("J107").Value = 7472.75
$ws.Range("K107").Value = 4128.7827
$ws.Range("L107").Value = 7472.75
$ws.Range("M107").Value = -2208.7827
$ws.Range("N107").Value = -11312.75

$ws.Range("H134").Value = 1976.2916
$ws.Range("I134").Value = 1720.2
$ws.Range("K134").Value = 5160.6
$ws.Range("M134").Value = -2625.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 161.55
$ws.Range("I7").Value = 17.545454
$ws.Range("J7").Value = 337.55554
$ws.Range("K7").Value = 17.545454
$ws.Range("L7").Value = 337.55554
$ws.Range("M7").Value = 95.45454599999999
$ws.Range("N7").Value = -563.5555400000001

$ws.Range("H22").Value = 431.27274
$ws.Range("I22").Value = 474.25
$ws.Range("K22").Value = 474.25
$ws.Range("M22").Value = -124.25

$ws.Range("H31").Value = 1516.5714
$ws.Range("J31").Value = 1950.375
$ws.Range("L31").Value = 1950.375
$ws.Range("N31").Value = -2540.375

$ws.Range("H34").Value = 1516.5714
$ws.Range("J34").Value = 1950.375
$ws.Range("L34").Value = 1950.375
$ws.Range("N34").Value = -2354.375

$ws.Range("H99").Value = 3182.6155
$ws.Range("I99").Value = 3028.7856
$ws.Range("K99").Value = 3028.7856
$ws.Range("M99").Value = -1530.7856

$ws.Range("H126").Value = 3182.6155
$ws.Range("I126").Value = 3028.7856
$ws.Range("K126").Value = 9086.356800000001
$ws.Range("M126").Value = -6616.356800000001

$ws.Range("H134").Value = 3302.1082
$ws.Range("I134").Value = 3340.3872
$ws.Range("K134").Value = 10021.1616
$ws.Range("M134").Value = -7486.161599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2390.818
$ws.Range("I68").Value = 959.8
$ws.Range("J68").Value = 3583.3333
$ws.Range("K68").Value = 2879.4
$ws.Range("L68").Value = 10749.9999
$ws.Range("M68").Value = -2068.4
$ws.Range("N68").Value = -12371.9999

$ws.Range("H71").Value = 2390.818
$ws.Range("I71").Value = 959.8
$ws.Range("J71").Value = 3583.3333
$ws.Range("K71").Value = 8638.199999999999
$ws.Range("L71").Value = 32249.9997
$ws.Range("M71").Value = -4582.199999999999
$ws.Range("N71").Value = -40361.9997

$ws.Range("H86").Value = 2777.8
$ws.Range("I86").Value = 2863
$ws.Range("J86").Value = 2650
$ws.Range("K86").Value = 8589
$ws.Range("L86").Value = 7950
$ws.Range("M86").Value = -7403
$ws.Range("N86").Value = -10322

$ws.Range("H89").Value = 2777.8
$ws.Range("I89").Value = 2863
$ws.Range("J89").Value = 2650
$ws.Range("K89").Value = 25767
$ws.Range("L89").Value = 23850
$ws.Range("M89").Value = -19839
$ws.Range("N89").Value = -35706

$ws.Range("H129").Value = 2105.6
$ws.Range("J129").Value = 2560
$ws.Range("L129").Value = 7680
$ws.Range("N129").Value = -17680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 230035
$ws.Range("I70").Value = 556805.5
$ws.Range("J70").Value = 12188
$ws.Range("K70").Value = 556805.5
$ws.Range("L70").Value = 12188
$ws.Range("M70").Value = -556535.5
$ws.Range("N70").Value = -12728

$ws.Range("H73").Value = 230035
$ws.Range("I73").Value = 556805.5
$ws.Range("J73").Value = 12188
$ws.Range("K73").Value = 556805.5
$ws.Range("L73").Value = 12188
$ws.Range("M73").Value = -555869.5
$ws.Range("N73").Value = -14060

$ws.Range("H122").Value = 6549.75
$ws.Range("I122").Value = 7974.75
$ws.Range("J122").Value = 5124.75
$ws.Range("K122").Value = 23924.25
$ws.Range("L122").Value = 15374.25
$ws.Range("M122").Value = -21474.25
$ws.Range("N122").Value = -20274.25

$ws.Range("H126").Value = 1122.1666
$ws.Range("I126").Value = 702.5
$ws.Range("J126").Value = 1961.5
$ws.Range("K126").Value = 2107.5
$ws.Range("L126").Value = 5884.5
$ws.Range("M126").Value = 362.5
$ws.Range("N126").Value = -10824.5

$ws.Range("H132").Value = 5979
$ws.Range("I132").Value = 5458.5576
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 16375.6728
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -13845.6728
$ws.Range("N132").Value = -50060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10102419
$ws.Range("I22").Value = 22727822
$ws.Range("K22").Value = 22727822
$ws.Range("M22").Value = -22727527

$ws.Range("H27").Value = 10102419
$ws.Range("I27").Value = 22727822
$ws.Range("K27").Value = 22727822
$ws.Range("M27").Value = -22727715

$ws.Range("H122").Value = 5479.6816
$ws.Range("I122").Value = 5104.9375
$ws.Range("J122").Value = 6479
$ws.Range("K122").Value = 15314.8125
$ws.Range("L122").Value = 19437
$ws.Range("M122").Value = -12864.8125
$ws.Range("N122").Value = -24337

$ws.Range("H139").Value = 38563.168
$ws.Range("J139").Value = 38563.168
$ws.Range("L139").Value = 38563.168
$ws.Range("N139").Value = -48843.168

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14651.75
$ws.Range("I45").Value = 7967
$ws.Range("J45").Value = 15606.714
$ws.Range("K45").Value = 7967
$ws.Range("L45").Value = 15606.714
$ws.Range("M45").Value = -7476
$ws.Range("N45").Value = -16588.714

$ws.Range("H52").Value = 19656.857
$ws.Range("I52").Value = 25999.5
$ws.Range("J52").Value = 17119.8
$ws.Range("K52").Value = 25999.5
$ws.Range("L52").Value = 17119.8
$ws.Range("M52").Value = -25773.5
$ws.Range("N52").Value = -17571.8

$ws.Range("H122").Value = 37038910
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
